$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the 13 "separator" paragraphs: a paragraph whose sole content is
#    a run of 60 U+2500 (BOX DRAWINGS LIGHT HORIZONTAL) characters, used as a
#    visual divider between sections.
#    We include the *leading* paragraph mark (the end of the PRECEDING
#    paragraph) in the search/replace pattern and replace with nothing; this
#    deletes the separator paragraph's own mark while leaving the following
#    paragraph's own pPr/style untouched (Word keeps the formatting that
#    belongs to the paragraph mark that survives the merge).
# ---------------------------------------------------------------------------
$sep = ""
for ($i = 0; $i -lt 60; $i++) {
    $sep += [char]0x2500
}
$sepPattern = "`r" + $sep
$findRange = $d.Content
$findRange.Find.Execute($sepPattern, $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Remove the 13 empty paragraphs that only carry <w:pPr><w:spacing
#    w:before="40"/></w:pPr> (no runs at all) which Word always placed right
#    after a table. They are identified by: paragraph text is just the
#    paragraph mark (length 1) and Format.SpaceBefore = 2 points (=40
#    twentieths of a point).
# ---------------------------------------------------------------------------
$emptyStarts = @()
$total = $d.Paragraphs.Count
for ($i = 1; $i -le $total; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Format.SpaceBefore -eq 2) {
        $t = $p.Range.Text
        if ($t.Length -eq 1) {
            $emptyStarts += $p.Range.Start
        }
    }
}
$emptyStarts = $emptyStarts | Sort-Object -Descending
foreach ($s in $emptyStarts) {
    $delRange = $d.Range($s - 1, $s + 1)
    $delRange.Delete()
}

# ---------------------------------------------------------------------------
# 3) Remove the 3 paragraphs that contain nothing but an inline picture.
# ---------------------------------------------------------------------------
$shapes = $d.InlineShapes
$imgStarts = @()
$shapeCount = $shapes.Count
for ($i = 1; $i -le $shapeCount; $i++) {
    $imgStarts += $shapes.Item($i).Range.Start
}
$imgStarts = $imgStarts | Sort-Object -Descending
foreach ($s in $imgStarts) {
    $delRange = $d.Range($s - 1, $s + 1)
    $delRange.Delete()
}
